$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that must be updated from
# 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (rows 2-494).
$ws.Range("C2:C494").Value = 45175
